$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Update the beverage names (shift content up one row: Lakkalikoori is dropped,
# and two new products replace the old Ipoh Coffee / Laughing Lumberjack Lager)
$ws.Range("A2").Value = "Chang"
$ws.Range("A3").Value = "Cote de Blaye"
$ws.Range("A4").Value = "Steeleye Stout"
$ws.Range("A5").Value = "Rhonbrau Klosterbier"

# Keep the failure note in sync with the renamed beverage on row 4
$ws.Range("D4").Value = "Unable to find beveraged named 'Steeleye Stout'"

# Add an empty, taller row 6 as part of the page-layout tweak
$ws.Rows.Item(6).RowHeight = 43.5
